# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '57.149.22'
Set-TextValue $ws.Range('E2') '  +1.38%  '
Set-TextValue $ws.Range('D3') '3.259.89'
Set-TextValue $ws.Range('E3') '  +0.77%  '
Set-TextValue $ws.Range('E4') '  +0.07%  '
Set-TextValue $ws.Range('D5') '397.75'
Set-TextValue $ws.Range('E5') '  -0.29%  '
Set-TextValue $ws.Range('D6') '108.82'
Set-TextValue $ws.Range('E6') '  -2.02%  '
Set-TextValue $ws.Range('E7') '  +4.23%  '
Set-TextValue $ws.Range('E8') '  +0.02%  '
Set-TextValue $ws.Range('E9') '  -0.80%  '
Set-TextValue $ws.Range('D10') '39.28'
Set-TextValue $ws.Range('E10') '  -0.58%  '
Set-TextValue $ws.Range('D11') '0.0953'
Set-TextValue $ws.Range('E11') '  +5.06%  '
Set-TextValue $ws.Range('E12') '  +1.56%  '
Set-TextValue $ws.Range('D13') '3.774.47'
Set-TextValue $ws.Range('E13') '  +0.96%  '
Set-TextValue $ws.Range('D14') '8.27'
Set-TextValue $ws.Range('E14') '  +1.84%  '
Set-TextValue $ws.Range('D15') '18.97'
Set-TextValue $ws.Range('E15') '  -0.65%  '
Set-TextValue $ws.Range('D16') '3.261.14'
Set-TextValue $ws.Range('E16') '  +2.73%  '
Set-TextValue $ws.Range('E17') '  -2.31%  '
Set-TextValue $ws.Range('D18') '11.08'
Set-TextValue $ws.Range('E18') '  +3.63%  '
Set-TextValue $ws.Range('D19') '56.936.40'
Set-TextValue $ws.Range('E19') '  +1.40%  '
Set-TextValue $ws.Range('E20') '  -1.20%  '
Set-TextValue $ws.Range('D21') '0.0000107'
Set-TextValue $ws.Range('E21') '  +5.07%  '
Set-TextValue $ws.Range('D22') '12.93'
Set-TextValue $ws.Range('E22') '  -1.10%  '
Set-TextValue $ws.Range('D23') '294.82'
Set-TextValue $ws.Range('D24') '73.93'
Set-TextValue $ws.Range('E24') '  -2.08%  '
Set-TextValue $ws.Range('E25') '  -1.60%  '
Set-TextValue $ws.Range('D26') '28.08'
Set-TextValue $ws.Range('E26') '  -0.79%  '
Set-TextValue $ws.Range('D27') '7.88'
Set-TextValue $ws.Range('E27') '  -3.95%  '
Set-TextValue $ws.Range('E28') '  +0.68%  '
Set-TextValue $ws.Range('D29') '7.46'
Set-TextValue $ws.Range('E29') '  -0.25%  '
Set-TextValue $ws.Range('E30') '  -2.85%  '
Set-TextValue $ws.Range('B31') 'Dai'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D31') '1.00'
Set-TextValue $ws.Range('E31') '  +0.04%  '
Set-TextValue $ws.Range('B32') 'Hedera'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D32') '0.113'
Set-TextValue $ws.Range('E32') '  +1.95%  '
Set-TextValue $ws.Range('D33') '11.18'
Set-TextValue $ws.Range('E33') '  -0.18%  '
Set-TextValue $ws.Range('D34') '40.24'
Set-TextValue $ws.Range('E34') '  +10.46%  '
Set-TextValue $ws.Range('D35') '0.0495'
Set-TextValue $ws.Range('E35') '  +0.41%  '
Set-TextValue $ws.Range('E36') '  +0.70%  '
Set-TextValue $ws.Range('D37') '51.31'
Set-TextValue $ws.Range('E37') '  -0.16%  '
Set-TextValue $ws.Range('E38') '  +0.01%  '
Set-TextValue $ws.Range('D39') '3.47'
Set-TextValue $ws.Range('E39') '  -1.65%  '
Set-TextValue $ws.Range('D40') '3.02'
Set-TextValue $ws.Range('E40') '  -3.99%  '
Set-TextValue $ws.Range('D41') '137.95'
Set-TextValue $ws.Range('E41') '  +2.32%  '
Set-TextValue $ws.Range('E42') '  +1.21%  '
Set-TextValue $ws.Range('D43') '0.285'
Set-TextValue $ws.Range('E43') '  +0.07%  '
Set-TextValue $ws.Range('E44') '  -2.95%  '
Set-TextValue $ws.Range('D45') '3.89'
Set-TextValue $ws.Range('E45') '  -3.76%  '
Set-TextValue $ws.Range('D46') '16.74'
Set-TextValue $ws.Range('E46') '  -2.88%  '
Set-TextValue $ws.Range('D47') '22.35'
Set-TextValue $ws.Range('E47') '  +0.30%  '
Set-TextValue $ws.Range('E48') '  +4.17%  '
Set-TextValue $ws.Range('B49') 'Maker'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D49') '2.147.52'
Set-TextValue $ws.Range('E49') '  +0.35%  '
Set-TextValue $ws.Range('B50') 'ApeXProtocol'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D50') '2.47'
Set-TextValue $ws.Range('E50') '  +0.18%  '
Set-TextValue $ws.Range('E51') '  -7.28%  '
